# Adds a new "IMAGENS:" section (with a pexels.com text entry) at the
# end of the document, right after the "ICONES:" section's
# "material.io (e do google)" paragraph and before the two trailing blank
# paragraphs, matching the pattern used by the other sections in the file.

$d = $word.ActiveDocument

# The anchor is the paragraph that ends the "ICONES" section (the one
# containing "material.io (e do google)") -- it is the paragraph
# immediately before the two trailing empty paragraphs at the end of the
# document.
$total = $d.Paragraphs.Count
$anchorIndex = $total - 2
$anchor = $d.Paragraphs.Item($anchorIndex)

$r = $anchor.Range
$r.Collapse(0)            # wdCollapseEnd

# 1) blank separator paragraph (inherits the sz24 formatting already on the
#    anchor paragraph's mark)               -> becomes paragraph anchorIndex+1
$r.InsertParagraphAfter()
$r.Collapse(0)

# 2) "IMAGENS:" heading paragraph (bold, sz26 -- same style as the other
#    section headings such as "ICONES:")     -> becomes paragraph anchorIndex+2
$r.InsertParagraphAfter()
$r.Collapse(0)

# 3) "pexels.com" paragraph (sz24 -- same style as the other link/text
#    paragraphs)                             -> becomes paragraph anchorIndex+3
$r.InsertParagraphAfter()
$r.Collapse(0)

# 4) trailing blank separator paragraph (sz24) -> becomes paragraph anchorIndex+4
$r.InsertParagraphAfter()
$r.Collapse(0)

# --- fill in the heading paragraph ------------------------------------
$headingIndex = $anchorIndex + 2
$heading = $d.Paragraphs.Item($headingIndex)
$hr = $heading.Range
$hStart = $hr.Start

$hr.Font.Bold = 1
$hr.Font.BoldBi = 1
$hr.Font.Size = 13
$hr.Font.SizeBi = 13
$hr.InsertAfter("IMAGENS")

$afterWord = $d.Range($hStart + 7, $hStart + 7)
$afterWord.Font.Bold = 1
$afterWord.Font.BoldBi = 1
$afterWord.Font.Size = 13
$afterWord.Font.SizeBi = 13
$afterWord.InsertAfter(":")

# --- fill in the "pexels.com" paragraph --------------------------------
$linkIndex = $anchorIndex + 3
$link = $d.Paragraphs.Item($linkIndex)
$lr = $link.Range
$lr.InsertAfter("pexels.com")

Write-Output "Done. Paragraphs now: $($d.Paragraphs.Count)"
